$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 currently holds WALA771012HCRGR054 / Wednesday Addams
# Row 5 currently holds GOTW771012HMRGR087 / Khal Drogo
# The new test record (GOTW.../Khal Drogo) is inserted right after "Walter White",
# which pushes WALA.../Wednesday Addams down to row 5 and gives it a purchase count of 1.

$ws.Range("A4").Value = "GOTW771012HMRGR087"
$ws.Range("B4").Value = "Khal Drogo"

$ws.Range("A5").Value = "WALA771012HCRGR054"
$ws.Range("B5").Value = "Wednesday Addams"
$ws.Range("C5").Value = 1
